$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_1")

# Camera config file name changed
$ws.Range("B6").Value = "Cam_Stereo.sdf"

# Marker pose value replaced by a plain numeric value
$ws.Range("C7").Value = 2

$ws.Select()
$ws.Range("B38").Select()
